$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old Row 18
$ws.Range("D2").Value2 = 44235
$ws.Range("L2").Value2 = "Primera"
$ws.Range("M2").Value2 = 60
$ws.Range("N2").Value2 = 3000
$ws.Range("O2").Value2 = 3000
$ws.Range("P2").Value2 = 3000
$ws.Range("Q2").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R2").Value2 = "Provincia de Linares"
$ws.Range("S2").Value2 = 1500
$ws.Range("T2").Value2 = 2

# Row 3 <- old Row 4
$ws.Range("D3").Value2 = 44204
$ws.Range("L3").Value2 = "Primera"
$ws.Range("M3").Value2 = 50
$ws.Range("N3").Value2 = 3000
$ws.Range("O3").Value2 = 3000
$ws.Range("P3").Value2 = 3000
$ws.Range("Q3").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R3").Value2 = "Provincia de Linares"
$ws.Range("S3").Value2 = 1500
$ws.Range("T3").Value2 = 2

# Row 4 <- old Row 5
$ws.Range("D4").Value2 = 44204
$ws.Range("L4").Value2 = "Segunda"
$ws.Range("M4").Value2 = 140
$ws.Range("N4").Value2 = 2400
$ws.Range("O4").Value2 = 2400
$ws.Range("P4").Value2 = 2400
$ws.Range("Q4").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R4").Value2 = "Provincia de Linares"
$ws.Range("S4").Value2 = 1200
$ws.Range("T4").Value2 = 2

# Row 5 <- old Row 17
$ws.Range("D5").Value2 = 44186
$ws.Range("L5").Value2 = "Primera"
$ws.Range("M5").Value2 = 200
$ws.Range("N5").Value2 = 3000
$ws.Range("O5").Value2 = 3000
$ws.Range("P5").Value2 = 3000
$ws.Range("Q5").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R5").Value2 = "Provincia de Limarí"
$ws.Range("S5").Value2 = 1500
$ws.Range("T5").Value2 = 2

# Row 6 <- old Row 9
$ws.Range("D6").Value2 = 44202
$ws.Range("L6").Value2 = "Primera"
$ws.Range("M6").Value2 = 30
$ws.Range("N6").Value2 = 3000
$ws.Range("O6").Value2 = 3000
$ws.Range("P6").Value2 = 3000
$ws.Range("Q6").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R6").Value2 = "Provincia de Linares"
$ws.Range("S6").Value2 = 1500
$ws.Range("T6").Value2 = 2

# Row 7 <- old Row 10
$ws.Range("D7").Value2 = 44202
$ws.Range("L7").Value2 = "Segunda"
$ws.Range("M7").Value2 = 20
$ws.Range("N7").Value2 = 2600
$ws.Range("O7").Value2 = 2600
$ws.Range("P7").Value2 = 2600
$ws.Range("Q7").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R7").Value2 = "Provincia de Linares"
$ws.Range("S7").Value2 = 1300
$ws.Range("T7").Value2 = 2

# Row 8 <- old Row 6
$ws.Range("D8").Value2 = 44232
$ws.Range("L8").Value2 = "Primera"
$ws.Range("M8").Value2 = 60
$ws.Range("N8").Value2 = 3000
$ws.Range("O8").Value2 = 3000
$ws.Range("P8").Value2 = 3000
$ws.Range("Q8").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R8").Value2 = "Provincia de Linares"
$ws.Range("S8").Value2 = 1500
$ws.Range("T8").Value2 = 2

# Row 9 <- old Row 12
$ws.Range("D9").Value2 = 44169
$ws.Range("L9").Value2 = "Primera"
$ws.Range("M9").Value2 = 400
$ws.Range("N9").Value2 = 3600
$ws.Range("O9").Value2 = 3600
$ws.Range("P9").Value2 = 3600
$ws.Range("Q9").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R9").Value2 = "Provincia de Linares"
$ws.Range("S9").Value2 = 1800
$ws.Range("T9").Value2 = 2

# Row 10 <- old Row 15
$ws.Range("D10").Value2 = 44200
$ws.Range("L10").Value2 = "Segunda"
$ws.Range("M10").Value2 = 50
$ws.Range("N10").Value2 = 2600
$ws.Range("O10").Value2 = 2600
$ws.Range("P10").Value2 = 2600
$ws.Range("Q10").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R10").Value2 = "Provincia de Linares"
$ws.Range("S10").Value2 = 1300
$ws.Range("T10").Value2 = 2

# Row 11 <- old Row 22
$ws.Range("D11").Value2 = 44167
$ws.Range("L11").Value2 = "Primera"
$ws.Range("M11").Value2 = 500
$ws.Range("N11").Value2 = 3600
$ws.Range("O11").Value2 = 3600
$ws.Range("P11").Value2 = 3600
$ws.Range("Q11").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R11").Value2 = "Región de O'Higgins"
$ws.Range("S11").Value2 = 1800
$ws.Range("T11").Value2 = 2

# Row 12 <- old Row 20
$ws.Range("D12").Value2 = 44210
$ws.Range("L12").Value2 = "Segunda"
$ws.Range("M12").Value2 = 150
$ws.Range("N12").Value2 = 2700
$ws.Range("O12").Value2 = 2700
$ws.Range("P12").Value2 = 2700
$ws.Range("Q12").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R12").Value2 = "Provincia de Linares"
$ws.Range("S12").Value2 = 1350
$ws.Range("T12").Value2 = 2

# Row 13 <- old Row 11
$ws.Range("D13").Value2 = 44165
$ws.Range("L13").Value2 = "Primera"
$ws.Range("M13").Value2 = 400
$ws.Range("N13").Value2 = 3400
$ws.Range("O13").Value2 = 3400
$ws.Range("P13").Value2 = 3400
$ws.Range("Q13").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R13").Value2 = "Región de O'Higgins"
$ws.Range("S13").Value2 = 1700
$ws.Range("T13").Value2 = 2

# Row 14 <- old Row 13
$ws.Range("D14").Value2 = 44172
$ws.Range("L14").Value2 = "Primera"
$ws.Range("M14").Value2 = 300
$ws.Range("N14").Value2 = 3400
$ws.Range("O14").Value2 = 3600
$ws.Range("P14").Value2 = 3467
$ws.Range("Q14").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R14").Value2 = "Provincia de Linares"
$ws.Range("S14").Value2 = 1734
$ws.Range("T14").Value2 = 2

# Row 15 <- old Row 16
$ws.Range("D15").Value2 = 44265
$ws.Range("L15").Value2 = "Primera"
$ws.Range("M15").Value2 = 70
$ws.Range("N15").Value2 = 3600
$ws.Range("O15").Value2 = 3800
$ws.Range("P15").Value2 = 3714
$ws.Range("Q15").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R15").Value2 = "Provincia de Linares"
$ws.Range("S15").Value2 = 1857
$ws.Range("T15").Value2 = 2

# Row 16 <- old Row 14
$ws.Range("D16").Value2 = 44187
$ws.Range("L16").Value2 = "Primera"
$ws.Range("M16").Value2 = 110
$ws.Range("N16").Value2 = 2600
$ws.Range("O16").Value2 = 3000
$ws.Range("P16").Value2 = 2782
$ws.Range("Q16").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R16").Value2 = "Provincia de Linares"
$ws.Range("S16").Value2 = 1391
$ws.Range("T16").Value2 = 2

# Row 17 <- old Row 19
$ws.Range("D17").Value2 = 44264
$ws.Range("L17").Value2 = "Primera"
$ws.Range("M17").Value2 = 110
$ws.Range("N17").Value2 = 3500
$ws.Range("O17").Value2 = 4000
$ws.Range("P17").Value2 = 3727
$ws.Range("Q17").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R17").Value2 = "Provincia de Linares"
$ws.Range("S17").Value2 = 1864
$ws.Range("T17").Value2 = 2

# Row 18 <- old Row 7
$ws.Range("D18").Value2 = 44166
$ws.Range("L18").Value2 = "Primera"
$ws.Range("M18").Value2 = 1500
$ws.Range("N18").Value2 = 3600
$ws.Range("O18").Value2 = 3600
$ws.Range("P18").Value2 = 3600
$ws.Range("Q18").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R18").Value2 = "Región de O'Higgins"
$ws.Range("S18").Value2 = 1800
$ws.Range("T18").Value2 = 2

# Row 19 <- old Row 8
$ws.Range("D19").Value2 = 44162
$ws.Range("L19").Value2 = "Primera"
$ws.Range("M19").Value2 = 100
$ws.Range("N19").Value2 = 4000
$ws.Range("O19").Value2 = 4000
$ws.Range("P19").Value2 = 4000
$ws.Range("Q19").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R19").Value2 = "Región de O'Higgins"
$ws.Range("S19").Value2 = 2000
$ws.Range("T19").Value2 = 2

# Row 20 <- old Row 21
$ws.Range("D20").Value2 = 44176
$ws.Range("L20").Value2 = "Primera"
$ws.Range("M20").Value2 = 150
$ws.Range("N20").Value2 = 3500
$ws.Range("O20").Value2 = 3500
$ws.Range("P20").Value2 = 3500
$ws.Range("Q20").Value2 = "`$/bandeja 12 canastillos 125 gramos"
$ws.Range("R20").Value2 = "Provincia de Curicó"
$ws.Range("S20").Value2 = 2333
$ws.Range("T20").Value2 = 1.5

# Row 21 <- old Row 2
$ws.Range("D21").Value2 = 44211
$ws.Range("L21").Value2 = "Primera"
$ws.Range("M21").Value2 = 40
$ws.Range("N21").Value2 = 2800
$ws.Range("O21").Value2 = 2800
$ws.Range("P21").Value2 = 2800
$ws.Range("Q21").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R21").Value2 = "Provincia de Linares"
$ws.Range("S21").Value2 = 1400
$ws.Range("T21").Value2 = 2

# Row 22 <- old Row 3
$ws.Range("D22").Value2 = 44211
$ws.Range("L22").Value2 = "Segunda"
$ws.Range("M22").Value2 = 30
$ws.Range("N22").Value2 = 2600
$ws.Range("O22").Value2 = 2600
$ws.Range("P22").Value2 = 2600
$ws.Range("Q22").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R22").Value2 = "Provincia de Linares"
$ws.Range("S22").Value2 = 1300
$ws.Range("T22").Value2 = 2
